$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the training schedule values on row 2 (x_corrSteps, y_nrSteps, alienID)
$ws.Range("D2").Value = 4
$ws.Range("F2").Value = -3
$ws.Range("H2").Value = 46

# Move the active selection to E2, matching the saved cursor position
$ws.Range("E2").Select()
